# Actualización automática 2025-06-03 09:15:07
$wb = $excel.ActiveWorkbook

$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# "VENTAS POR GRUPO": K3 gets the sale amount, K7 progress label bumps to 1 de 5
$wsGrupo.Range("K3").Value = 851.4299999999999
$wsGrupo.Range("K7").Value = "1 de 5"

# "VENTA MENSUAL": column F widens slightly and picks up the same sale amount
$wsMensual.Range("F3").Value = 851.4299999999999
$wsMensual.Range("F7").Value = 851.4299999999999
# Excel's ColumnWidth (chars) stores as width+5/6 px-padding in the OOXML
# <col> width attribute, so back the padding out to land on width=12 exactly.
$wsMensual.Columns.Item(6).ColumnWidth = 12 - (5 / 6)
